$wb = $excel.ActiveWorkbook
$wsQ = $wb.Worksheets("Preguntas")
$wsC = $wb.Worksheets("Conocimiento")

# Read the 10 question texts from Preguntas!B2:B11 (preserves special chars like nbsp)
$questions = @()
for ($r = 2; $r -le 11; $r++) {
    $questions += $wsQ.Cells.Item($r, 2).Value()
}

# Overwrite Conocimiento header row A1:J1 with the long question texts
for ($c = 1; $c -le 10; $c++) {
    $wsC.Cells.Item(1, $c).Value = $questions[$c - 1]
}

# Set K1 header to "Decision"
$wsC.Cells.Item(1, 11).Value = "Decisión"

# Fix J5 value from 1 to 0
$wsC.Cells.Item(5, 10).Value = 0

# Update selection to K1
$wsC.Range("K1").Select()

# Remove the Preguntas sheet entirely
$wsQ.Delete()

# Re-fetch the surviving sheet (reference can go stale after a sibling delete)
# and rename it to Hoja1
$wsC2 = $wb.Worksheets("Conocimiento")
$wsC2.Name = "Hoja1"

Write-Output "done"
